# Recalculate the mean / refresh the accuracy values for the GlobalThresholding
# (HK_G_acc_LT) worksheet. This re-writes the header label (which causes Excel
# to re-intern the shared string for A1) and refreshes the handful of subject
# accuracy values whose underlying counts changed after the mean was
# recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the header text so the cell's shared-string reference is rewritten.
$ws.Range("A1").Value = "HK_G_acc_LT"

# Updated accuracy figures (recomputed means) for the affected subjects.
$ws.Range("A34").Value = 74.831081081081081
$ws.Range("A39").Value = 65.03378378378379
$ws.Range("A49").Value = 74.155405405405403
